# Weekly update: insert a new "Ajo" (garlic) price record for
# Agrícola del Norte S.A. de Arica as the new row 23, pushing the
# existing rows 23-38 down to 24-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 23 (shifts 23..38 -> 24..39,
# and copies formatting - e.g. the date number format on column D -
# down from the row above, same as Excel's native "Insert" behaviour).
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with this week's entry.
$ws.Cells.Item(23, 1).Value = 1
$ws.Cells.Item(23, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(23, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(23, 4).Value = 45001
$ws.Cells.Item(23, 5).Value = 15
$ws.Cells.Item(23, 6).Value = 100112003
$ws.Cells.Item(23, 7).Value = "Ajo"
$ws.Cells.Item(23, 8).Value = "Chino"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 400
$ws.Cells.Item(23, 11).Value = 18000
$ws.Cells.Item(23, 12).Value = 20000
$ws.Cells.Item(23, 13).Value = 19000
$ws.Cells.Item(23, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(23, 15).Value = "China"
$ws.Cells.Item(23, 16).Value = 1900
$ws.Cells.Item(23, 17).Value = 10
$ws.Cells.Item(23, 18).Value = "Hortaliza"
